# CanavanCalculator2 - "better error handling, data validity checking"
#
# The "Bazinga" scratch/demo sheet is replaced as the example data source by a
# new "Test Player" sheet reference on the Calculations sheet, and the README
# sheet gets a clarifying note about right field numbering next to the
# outcome-code legend.

$wb = $excel.ActiveWorkbook

$wsReadme = $wb.Worksheets.Item("README")
$wsCalc   = $wb.Worksheets.Item("Calculations")

# --- Calculations!E8: point the calculator at the "Test Player" example sheet
# instead of the old "Bazinga" scratch sheet (drives every INDIRECT() lookup
# on the Calculations sheet, so all the cached stats below recompute too).
$wsCalc.Range("E8").Value = "Test Player"

# --- README: add a clarifying note beside the outcome-code table explaining
# that code 9 really means right field.
$wsReadme.Range("L21:N21").Merge()
$wsReadme.Range("L21").Value = "<--- It's really Right Field for 9"
$wsReadme.Range("L21:N21").HorizontalAlignment = -4108

# --- Restore the on-screen selections / active sheet left by the author.
$wsCalc.Activate()
$wsCalc.Range("E9").Select()

$wsReadme.Activate()
$wsReadme.Range("M24").Select()
